$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.344.48"
$ws.Range("E2").Value = "  +1.10%  "

# Row 3
$ws.Range("D3").Value = "1.667.18"
$ws.Range("E3").Value = "  +0.94%  "

# Row 4
$ws.Range("E4").Value = "  +0.93%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.52"
$ws.Range("E5").Value = "  +1.04%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5347"
$ws.Range("E6").Value = "  +1.32%  "

# Row 7
$ws.Range("E7").Value = "  +0.85%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2666"
$ws.Range("E8").Value = "  +2.66%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06396"
$ws.Range("E9").Value = "  +1.33%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.91"
$ws.Range("E10").Value = "  +3.00%  "

# Row 11
$ws.Range("E11").Value = "  +0.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.561"
$ws.Range("E12").Value = "  +1.01%  "

# Row 13
$ws.Range("D13").Value = "1.666.93"
$ws.Range("E13").Value = "  +0.91%  "

# Row 14
$ws.Range("D14").Value = "1.895.78"
$ws.Range("E14").Value = "  +0.92%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5543"
$ws.Range("E15").Value = "  +1.26%  "

# Row 16
$ws.Range("D16").Value = "0.0₅8191"
$ws.Range("E16").Value = "  +0.15%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.02"
$ws.Range("E17").Value = "  +0.92%  "

# Row 18
$ws.Range("D18").Value = "26.375.61"
$ws.Range("E18").Value = "  +1.17%  "

# Row 19
$ws.Range("E19").Value = "  +0.89%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.678"
$ws.Range("E20").Value = "  +2.32%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.55"
$ws.Range("E21").Value = "  +2.58%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.28"
$ws.Range("E22").Value = "  +2.16%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.044"
$ws.Range("E23").Value = "  +0.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.011"
$ws.Range("E24").Value = "  +0.92%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.03"
$ws.Range("E25").Value = "  +1.78%  "

# Row 26
$ws.Range("E26").Value = "  -0.42%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.227"
$ws.Range("E27").Value = "  +0.24%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.18"
$ws.Range("E28").Value = "  +1.12%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.499"
$ws.Range("E29").Value = "  +4.58%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05868"
$ws.Range("E30").Value = "  +1.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.285"
$ws.Range("E31").Value = "  +1.21%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.584"
$ws.Range("E32").Value = "  +1.23%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.291"
$ws.Range("E33").Value = "  +1.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.613"
$ws.Range("E34").Value = "  +1.36%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9718"
$ws.Range("E35").Value = "  +3.20%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.832"
$ws.Range("E36").Value = "  +1.33%  "

# Row 37
$ws.Range("E37").Value = "  +0.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5829"
$ws.Range("E38").Value = "  +1.44%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01610"
$ws.Range("E39").Value = "  +0.40%  "

# Row 40
$ws.Range("D40").Value = "1.077.49"
$ws.Range("E40").Value = "  +4.97%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8644"
$ws.Range("E41").Value = "  +1.94%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.854"
$ws.Range("E42").Value = "  +2.47%  "

# Row 43
$ws.Range("E43").Value = "  +0.90%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.32"
$ws.Range("E44").Value = "  -0.23%  "

# Row 45
$ws.Range("D45").Value = "1.805.68"
$ws.Range("E45").Value = "  +0.62%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.10"
$ws.Range("E46").Value = "  +1.69%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.015"
$ws.Range("E47").Value = "  +1.39%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4392"
$ws.Range("E48").Value = "  +1.47%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.029"
$ws.Range("E49").Value = "  +2.35%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₈101"
$ws.Range("E50").Value = "  -8.84%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05166"
$ws.Range("E51").Value = "  +0.54%  "

